$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("利欧股份", "利欧股份", "利欧股份"),
    @("浙文互联", "特变电工", "航天发展"),
    @("特变电工", "白银有色", "白银有色"),
    @("白银有色", "浙文互联", "浙文互联"),
    @("通鼎互联", "湖南黄金", "蓝色光标"),
    @("湖南黄金", "中国西电", "中超控股"),
    @("蓝色光标", "红 宝 丽", "湖南黄金"),
    @("西部材料", "保变电气", "红宝丽"),
    @("红 宝 丽", "西部材料", "湖南白银"),
    @("航天发展", "中超控股", "紫金矿业"),
    @("中超控股", "蓝色光标", "特变电工"),
    @("中国西电", "航天发展", "保变电气"),
    @("保变电气", "紫金矿业", "茂业商业"),
    @("紫金矿业", "贵州茅台", "西部材料"),
    @("湖南白银", "通鼎互联", "锋龙股份"),
    @("锋龙股份", "湖南白银", "通鼎互联"),
    @("皇台酒业", "锋龙股份", "皇台酒业"),
    @("网宿科技", "皇台酒业", "中国西电"),
    @("润泽科技", "双杰电气", "横店影视"),
    @("铜陵有色", "铜陵有色", "铜陵有色")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
